# Integration with track model for lights and crossing
# The "junctionNS" indicator column (D) moves from its default red state to
# either a green state or is cleared out, depending on the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Green_info")

# Rows whose junctionNS indicator flips from red to green.
$greenRows = @(2, 3, 78, 102, 152)
foreach ($r in $greenRows) {
    $ws.Range("D$r").Value = "junctionNS;green"
}

# Rows whose junctionNS indicator is removed entirely.
$clearRows = @(14, 32, 64, 88, 103)
foreach ($r in $clearRows) {
    $ws.Range("D$r").ClearContents()
}

# Reflect the editor's final scroll position / selection in the sheet view.
$win = $excel.ActiveWindow
$ws.Range("A71").Select()
$win.ScrollRow = 71
$win.ScrollColumn = 1
$ws.Range("D88").Select()
